$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DBD")

# Insert a new row above row 9 (old header row becomes row 10, etc.)
$ws.Range("A9:G9").Insert(-4121)

# New row 9 picks up the same look as row 8 (border/fill/font for the
# label column A/B and the value column C..G)
$ws.Range("A8:G8").Copy()
$ws.Range("A9:G9").PasteSpecial(-4122)

# Fill in the new index-definition row (string added first so the
# shared-string table keeps the same ordering as the source edit)
$ws.Range("C9").Value = "AcDate,TitaTlrNo,TitaTxtNo"
$ws.Range("A9").Value = "Index5"

# A9:B9 is merged just like the other label cells above it
$ws.Range("A9:B9").Merge()

# The ProcCode hyperlink cell moved from G27 down to G28 with the row shift
$ws.Range("G27").Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("G28"), "", "ProcCode!A1", [System.Reflection.Missing]::Value, "ProcCode!A1")

# Reset the view: scroll back to the top and select C8
$ws.Range("C8").Select() | Out-Null
